# Update the "Ganhador 1" (C), "Ganhador 2" (D) and "Empates" (E) columns
# for every data row (rows 2-101) with the new simulation results:
#   Ganhador 1: 59 -> 63
#   Ganhador 2: 29 -> 23
#   Empates   : 12 -> 14
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 101; $r++) {
    $ws.Cells.Item($r, 3).Value = 63
    $ws.Cells.Item($r, 4).Value = 23
    $ws.Cells.Item($r, 5).Value = 14
}
